# Update cryptos list: refresh Price (column D) and Volume(1h) (column E)
# values for rows 2-51 to match the latest scrape, per the
# "Updated cryptos list" GitHub Actions commit.
#
# Price cells that look numeric (a single decimal point) get forced to
# Text via NumberFormat "@" before the write, then the style is reset
# back to Normal so the cell keeps its original (default) style index -
# only the stored text content changes, matching the source XML diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.364.18'
$ws.Range("E2").Value = '  -2.43%  '
$ws.Range("D3").Value = '1.775.73'
$ws.Range("E3").Value = '  -1.00%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.61%  '
$ws.Range("E5").Value = '  -0.58%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '304.89'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4237'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.50%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3598'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.15%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07160'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.64%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8363'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.08%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.44'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.87%  '
$ws.Range("D12").Value = '1.765.92'
$ws.Range("E12").Value = '  -3.66%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.438'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.56%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.242'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.23%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06901'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.87%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.86%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '79.07'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.78%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008649'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.90%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.90'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.06%  '
$ws.Range("D21").Value = '26.371.79'
$ws.Range("E21").Value = '  -3.49%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.075'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.50%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.89'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.79%  '
$ws.Range("D24").Value = '1.987.50'
$ws.Range("E24").Value = '  -4.00%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.45'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.83%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.801'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -8.45%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.96'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.32%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.078'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.37%  '
$ws.Range("E29").Value = '  +1.65%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.831'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +10.52%  '
$ws.Range("E31").Value = '  -0.86%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7272'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.28%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.121'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.91%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.319'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.76%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.001'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.65%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.734'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.01%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.087'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.28%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05109'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.17%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01881'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.50%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4920'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.60%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1605'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.89%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.595'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.29%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.319'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.66%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.019'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.33%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '104.33'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.39%  '
$ws.Range("E46").Value = '  -0.61%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.14'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.18%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.630'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.57%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06168'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4444'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.01%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.720'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.51%  '
